$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), styled like the existing
# header row (copy format from H1, the last existing header cell).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2 and J2 (plain numeric values, no special style)
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6
